$d = $word.ActiveDocument

$xml18 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">    "We are steadily adding cities and analyses to this resource with the aim of achieving better models that can support better depression public health policy and treatment practices that anticipate the likely dynamic evolution of the disease."</w:t></w:r></w:p>
'@
$d.Paragraphs(18).Range.InsertXML($xml18) | Out-Null

$xml16 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">    "We have incorporated interactive plots to allow users to explore the data and findings in a more engaging and informative way. These plots enable users to visualize depression trends across different LSOAs, compare the situation and results between cities, and understand the predictive </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>usefuleness</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> of different factors in explaining depression prevalence and growth.",</w:t></w:r></w:p>
'@
$d.Paragraphs(16).Range.InsertXML($xml16) | Out-Null

$xml14 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">    "This outcome falls short of our goal to provide a robust basis for predicting Depression Growth Events.  However, it is an intriguing result that is guiding us to investigate further using more sophisticated models that relaxes the assumption of linearity and based on a more accurate description and theory of how depression works.  The theory is itself evolving as we observe outcomes of our experiments with </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>data analysis and modelling but generally based on Depression Epidemiology, Etiology and Prognosis.",</w:t></w:r></w:p>
'@
$d.Paragraphs(14).Range.InsertXML($xml14) | Out-Null

$xml12 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">    "Growth events have, as measured by R-squared error, produced by Geographically Weighted Regression (GWR) have a high positive correlation with the Prevelance of Depression, and Depression Prevelance has a strong inverse correlation with the number of Prescription Items per patient.  However, only a handful of the 4835 LSOA’s in Greater London also have corresponding t-values that would indicate statistical significance.",</w:t></w:r></w:p>
'@
$d.Paragraphs(12).Range.InsertXML($xml12) | Out-Null

$xml10 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">    "From these datasets we have found in Greater London that Depression and Prescription Prevelance appears to grow steadily and form clear sustained clusters of both high and low levels.  But growth events do not form obvious clusters or follow any obvious pattern… and yet the clusters that appear to form for Prevelance are outcomes from the same data as the growth.  This is an unresolved puzzle, why does growth appear to erupt randomly whereas the consequent prevalence appears orderly and clustered.",</w:t></w:r></w:p>
'@
$d.Paragraphs(10).Range.InsertXML($xml10) | Out-Null

$xml8 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">    "This analysis uses spatial data of Depression Cases and Prescriptions of anti-depressants within Lower Super Output Areas ‘LSOAs’ in England.  An LSOA is an  administrative area with 400 – 1,000 households and between 1,000 and 3,000 persons.",    </w:t></w:r></w:p>
'@
$d.Paragraphs(8).Range.InsertXML($xml8) | Out-Null

$xml6 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="709" w:hanging="851"/></w:pPr><w:r><w:t xml:space="preserve">        "3.  Research into the causes of depression cases has identified changes in circumstances as the dominant driver of new depression cases. Changes in circumstances within a geographic location are, by definition, are more likely to lead to growth events than steady state ambient conditions.", </w:t></w:r></w:p>
'@
$d.Paragraphs(6).Range.InsertXML($xml6) | Out-Null

$xml5 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="709" w:hanging="851"/></w:pPr><w:r><w:t xml:space="preserve">        "2.  To date all Spatial analysis investigations of Depression Prevelance we have found focus on a search for factors that prove a clear link between steady state ambient conditions and Depression Prevalence",</w:t></w:r></w:p>
'@
$d.Paragraphs(5).Range.InsertXML($xml5) | Out-Null

$xml4 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="709" w:hanging="851"/></w:pPr><w:r><w:t xml:space="preserve">        "1.  Assessing the explanatory power of drivers that may explain spatial variations of Depression Growth events addresses a major gap in the literature",</w:t></w:r></w:p>
'@
$d.Paragraphs(4).Range.InsertXML($xml4) | Out-Null

$xml3 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">    "This study is focussed on improving our understanding by revealing Depression Growth events as opposed to studying factors that might correlate with Depression Prevelance.  There are a number of reasons we have chosen this focus but chief among them:",</w:t></w:r></w:p>
'@
$d.Paragraphs(3).Range.InsertXML($xml3) | Out-Null

$xml1 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">    </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">    "Depression has evolved in cities through a complex dynamic process that in some cases cities have similarities and in other cases the developments are distinctive.  We are seeking to identify patterns in the spatial variations of depression evolution that may offer clues about the underlying complex dynamic processes.   Ultimately the goal is to  discover better explanations of the variety of developments that may lead to a more robust theory and predictive capability.",</w:t></w:r></w:p>
'@
$d.Paragraphs(1).Range.InsertXML($xml1) | Out-Null
